$d = $word.ActiveDocument

# Each of the three multiple-choice lists below has one answer option whose
# leading letter ("C.", "B.", "D.") should be highlighted in red to mark it
# as the correct answer, while the rest of the paragraph text is untouched.

# "C. Pattern - is a solution to a specific type of problems..."
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Text = "C. Pattern - is a solution"
$rng1.Find.Forward = $true
$rng1.Find.Wrap = 0
if ($rng1.Find.Execute()) {
    $letter1 = $rng1.Duplicate
    $letter1.SetRange($rng1.Start, $rng1.Start + 1)
    $letter1.Font.Color = 255
}

# "B. Software design pattern - is set of patterns which solve commonly occurring problem within a given context..."
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Text = "B. Software design pattern - is set of patterns which solve commonly occurring problem within"
$rng2.Find.Forward = $true
$rng2.Find.Wrap = 0
if ($rng2.Find.Execute()) {
    $letter2 = $rng2.Duplicate
    $letter2.SetRange($rng2.Start, $rng2.Start + 1)
    $letter2.Font.Color = 255
}

# "D. All of the above."
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Text = "D. All of the above."
$rng3.Find.Forward = $true
$rng3.Find.Wrap = 0
if ($rng3.Find.Execute()) {
    $letter3 = $rng3.Duplicate
    $letter3.SetRange($rng3.Start, $rng3.Start + 1)
    $letter3.Font.Color = 255
}
